$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("LF"/old formula column shifts right
# along with everything else). This becomes the new "Type" column.
$ws.Columns("B").Insert()

# Fix header typo and add the new "Type" header.
$ws.Range("A1").Value = "Heuristic"
$ws.Range("B1").Value = "Type"

# Populate the new "Type" column for each heuristic row.
$ws.Range("B2").Value = "continuous"
$ws.Range("B3").Value = "continuous"
$ws.Range("B4").Value = "continuous"
$ws.Range("B5").Value = "continuous"
$ws.Range("B6").Value = "continuous"
$ws.Range("B7").Value = "continuous"
$ws.Range("B8").Value = "continuous"
$ws.Range("B9").Value = "continuous"
$ws.Range("B10").Value = "categorical & continuous"
$ws.Range("B11").Value = "continuous"
$ws.Range("B12").Value = "continuous"

# Rewrite the "LF" condition text (now column C) into the new IF/FRAUD/Abstain wording.
$ws.Range("C2").Value = "IF (txn_volume_vs_income > 2477.8) Fraud, ELSE Abstain"
$ws.Range("C3").Value = "IF (txn_volume_vs_occupation_median > 1562.0) FRAUD, ELSE Abstain"
$ws.Range("C4").Value = "IF (median_amt_wire > 6817.2 & (age <= 25 | income <= 40000)) FRAUD, ELSE Abstain"
$ws.Range("C5").Value = "IF (wire_ratio > 0.47 & (age <= 25 | income <= 40000)) FRAUD, ELSE Abstain"
$ws.Range("C6").Value = "IF (count_txn_below_threshold_frequency > 0.009 & n_txn_total > 5) FRAUD, ELSE Abstain"
$ws.Range("C7").Value = "IF (n_txn_total > 50 & median_hold_time_funds <= 0.11) FRAUD, ELSE Abstain"
$ws.Range("C8").Value = "IF (cross_border_ratio > 0.17 & transaction_unique_countries >= 3) FRAUD, ELSE Abstain"
$ws.Range("C9").Value = "IF (transaction_ecommerce_ratio > 0.98 & transaction_volume_90d > 198059.4 & n_txn_total > 10) FRAUD, ELSE Abstain"
$ws.Range("C10").Value = "IF (transaction_ecommerce_ratio > 0.98 & occupation == UNEMPLOYED) FRAUD, ELSE Abstain"
$ws.Range("C11").Value = "IF (transaction_unique_merchants > 46 & n_txn_total / transaction_unique_merchants < 1.5) FRAUD, ELSE Abstain"
$ws.Range("C12").Value = "IF (transaction_same_amount_frequency_7d > 0.14 & transaction_round_amount_frequency_7d > 0.5) FRAUD, ELSE Abstain"

# Adjust column widths: new "Type" column narrower, and the condition column
# (old "Note" column) a bit wider than before.
$ws.Columns("B").ColumnWidth = 21.833333333333336
$ws.Columns("C").ColumnWidth = 50.83333333333333
